$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("_settings")

# A row insert shifts every row below the insertion point (including any
# "scratch" row we might stage data in), so stage the formats we need to
# reuse for the new row in an unused block of columns instead - those are
# unaffected by the row insert.
$ws.Range("A14").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("B9:D9").Copy()
$ws.Range("AA1:AC1").PasteSpecial(-4122)
$ws.Range("E10:G10").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Insert a new row above row 9 (table paramTable row), shifting rows 9-15 down to 10-16.
$ws.Rows.Item(9).Insert()

# Copy the staged formatting onto the newly inserted (now blank) row 9.
$ws.Range("Z1").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("AA1:AC1").Copy()
$ws.Range("B9:D9").PasteSpecial(-4122)
$ws.Range("AD1:AF1").Copy()
$ws.Range("E9:G9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Clean up the staging area.
$ws.Range("Z1:AF1").Clear()

# Populate the new row with the "posId" parameter entry.
$ws.Range("A9").Value = "posId"
$ws.Range("B9").Value = "posId"
$ws.Range("D9").Value = "provided"

# Expand the table / autofilter range to include the new row.
$tbl = $ws.ListObjects.Item("paramTable")
$tbl.Resize($ws.Range("A8:G16"))

Write-Output "done"
